{"js": "// \"Version 2.\" -> \"Version 1.\"\n// The paragraph starts as four runs: \"Versi\" | \"on\" | \" 2\" | \".\" with a\n// _GoBack bookmark sitting between the \" 2\" run and the trailing \".\" run.\n// We need to end up with two runs: \"Version\" | \" 1.\" and the bookmark\n// kept in place (immediately after the \" 1.\" run).\n\nconst body = context.document.body;\n\n// Step 1: \"Versi\" + \"on\" -> \"Version\".\n// search() matches across the run boundary; replacing the whole match\n// merges the touched runs into a single run containing \"Version\".\nlet found = body.search(\"Version\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nfound.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// Step 2: \" 2\" -> \" 1.\" (this match ends right before the bookmark, so\n// the bookmark between this run and the trailing \".\" run is untouched).\nfound = body.search(\" 2\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nfound.items[0].insertText(\" 1.\", \"Replace\");\nawait context.sync();\n\n// Step 3: delete the now-redundant trailing \".\" run (the paragraph's old\n// closing period, now duplicated since step 2 appended \".\" already).\n// There are two \".\" matches left (\"...1.\" and the lone trailing run) \u2014\n// take the last one so the bookmark (which sits right before it) stays.\nfound = body.search(\".\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nfound.items[found.items.length - 1].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Step 1: \"Versi\" + \"on\" -> single run \"Version\"\n# (Find matching the full word across the run boundary merges the two\n#  runs touched by the match into one, dropping the split.)\n$find1 = $d.Content.Find\n$find1.Text = \"Version\"\n$find1.Replacement.Text = \"Version\"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# Step 2: \" 2\" -> \" 1.\" (the match stops before the bookmark, so the\n# bookmark between this run and the trailing \".\" run is left in place)\n$find2 = $d.Content.Find\n$find2.Text = \" 2\"\n$find2.Replacement.Text = \" 1.\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# Step 3: remove the now-redundant trailing \".\" run (the old closing\n# period run) by deleting its single character range directly, so the\n# bookmark immediately before it is left untouched.\n$para = $d.Paragraphs(1).Range\n$tail = $d.Range($para.End - 2, $para.End - 1)\nif ($tail.Text -eq \".\") {\n    $tail.Delete()\n}\n"}
